# StructureDefinition-employee-wage-amount.xlsx update
# - Metadata sheet: URL, Version, Date, Publisher values refreshed for the
#   LinuxForHealth (formerly Alvearie) rebrand / 8.0.0 release.
# - Elements sheet: clear the stale combined ele-1/ext-1 "Constraint(s)"
#   text that had been duplicated onto the root "Extension" row.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-wage-amount"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the root "Extension" element; column AI is "Constraint(s)".
$elements.Range("AI2").Value = ""
